# Apply "custom accuracy" rounding to row 5 (B5:AH5) and remove the now
# redundant row 6 (data reduced from 1000 rows worth down by one sample row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the measurement values in row 5 to 2 decimal places (custom accuracy).
$ws.Range("B5").Value = 16.81
$ws.Range("C5").Value = 12.29
$ws.Range("D5").Value = 1.12
$ws.Range("E5").Value = 36.54
$ws.Range("F5").Value = 29.68
$ws.Range("G5").Value = 13.23
$ws.Range("H5").Value = 51.11
$ws.Range("I5").Value = 20.36
$ws.Range("J5").Value = 8.98
$ws.Range("K5").Value = 13.24
$ws.Range("L5").Value = 14.66
$ws.Range("M5").Value = 15.44
$ws.Range("N5").Value = 4.23
$ws.Range("O5").Value = 13.16
$ws.Range("P5").Value = 18.68
$ws.Range("Q5").Value = 11.18
$ws.Range("R5").Value = 0.81
$ws.Range("S5").Value = 0.73
$ws.Range("T5").Value = 192.97
$ws.Range("U5").Value = 36.8
$ws.Range("V5").Value = 12.15
$ws.Range("W5").Value = 24.64
$ws.Range("X5").Value = 12.9
$ws.Range("Y5").Value = 2.06
$ws.Range("Z5").Value = 24.85
$ws.Range("AA5").Value = 10.73
$ws.Range("AB5").Value = 9.56
$ws.Range("AC5").Value = 11.24
$ws.Range("AD5").Value = 15.37
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 46.48
$ws.Range("AG5").Value = 6.79
$ws.Range("AH5").Value = 15.19

# Remove the last data row (row 6); the sheet's used range shrinks to A1:AH5.
$ws.Rows.Item(6).Delete()
